# Updates cryptos list: refresh price/volume figures and re-sort swapped rows (Aug 13 2024 run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as plain text, preserving the cell's original (default) style
# so no spurious number formats / styles are introduced (data stays inline/shared string).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "58.789.00"
Set-TextValue "E2" "  +1.00%  "

Set-TextValue "D3" "2.629.01"
Set-TextValue "E3" "  +2.12%  "

Set-TextValue "E4" "  -0.01%  "

Set-TextValue "D5" "520.77"
Set-TextValue "E5" "  +2.66%  "

Set-TextValue "D6" "144.54"
Set-TextValue "E6" "  +0.12%  "

Set-TextValue "E7" "  +0.00%  "

Set-TextValue "D8" "0.569"
Set-TextValue "E8" "  -0.38%  "

Set-TextValue "D9" "2.636.21"
Set-TextValue "E9" "  +1.99%  "

Set-TextValue "D10" "6.30"
Set-TextValue "E10" "  +0.24%  "

Set-TextValue "E11" "  +0.86%  "

Set-TextValue "E12" "  -0.45%  "

Set-TextValue "E13" "  -0.98%  "

Set-TextValue "D14" "3.089.02"
Set-TextValue "E14" "  +2.09%  "

Set-TextValue "D15" "58.812.63"
Set-TextValue "E15" "  +1.06%  "

Set-TextValue "D16" "20.77"
Set-TextValue "E16" "  -1.20%  "

Set-TextValue "E17" "  -0.58%  "

Set-TextValue "D18" "2.631.85"
Set-TextValue "E18" "  +1.65%  "

Set-TextValue "D19" "346.26"
Set-TextValue "E19" "  +1.26%  "

Set-TextValue "E20" "  -2.09%  "

Set-TextValue "D21" "10.18"
Set-TextValue "E21" "  -1.14%  "

Set-TextValue "D22" "6.14"
Set-TextValue "E22" "  +1.23%  "

Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  +0.00%  "

Set-TextValue "D24" "61.38"
Set-TextValue "E24" "  +1.29%  "

Set-TextValue "D25" "0.415"
Set-TextValue "E25" "  -0.95%  "

Set-TextValue "B26" "Kaspa"
Set-TextValue "C26" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D26" "0.163"
Set-TextValue "E26" "  +2.75%  "

Set-TextValue "B27" "Binance-PegBSC-USD"
Set-TextValue "C27" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D27" "0.996"
Set-TextValue "E27" "  -0.15%  "

Set-TextValue "D28" "0.0₃0799"
Set-TextValue "E28" "  -1.89%  "

Set-TextValue "D29" "7.07"
Set-TextValue "E29" "  +1.03%  "

Set-TextValue "E30" "  +0.00%  "

Set-TextValue "D31" "6.22"
Set-TextValue "E31" "  +1.91%  "

Set-TextValue "D32" "18.84"
Set-TextValue "E32" "  +0.17%  "

Set-TextValue "E33" "  +2.22%  "

Set-TextValue "D34" "150.20"
Set-TextValue "E34" "  +0.76%  "

Set-TextValue "D35" "0.980"
Set-TextValue "E35" "  +3.12%  "

Set-TextValue "D36" "3.97"
Set-TextValue "E36" "  -0.29%  "

Set-TextValue "E37" "  +0.18%  "

Set-TextValue "D38" "36.66"
Set-TextValue "E38" "  +1.84%  "

Set-TextValue "D39" "0.837"
Set-TextValue "E39" "  -1.55%  "

Set-TextValue "B40" "Stacks"
Set-TextValue "C40" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D40" "1.42"
Set-TextValue "E40" "  +1.81%  "

Set-TextValue "B41" "Filecoin"
Set-TextValue "C41" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D41" "3.64"
Set-TextValue "E41" "  +1.58%  "

Set-TextValue "B42" "Bittensor"
Set-TextValue "C42" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D42" "279.38"
Set-TextValue "E42" "  -3.23%  "

Set-TextValue "B43" "FirstDigitalUSD"
Set-TextValue "C43" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "0.997"
Set-TextValue "E43" "  +0.17%  "

Set-TextValue "D44" "0.0981"
Set-TextValue "E44" "  -0.85%  "

Set-TextValue "D45" "0.599"
Set-TextValue "E45" "  -1.65%  "

Set-TextValue "D46" "19.51"
Set-TextValue "E46" "  +1.19%  "

Set-TextValue "B47" "Hedera"
Set-TextValue "C47" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D47" "0.0520"
Set-TextValue "E47" "  -2.83%  "

Set-TextValue "B48" "WhiteBITCoin"
Set-TextValue "C48" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D48" "10.31"
Set-TextValue "E48" "  +0.61%  "

Set-TextValue "D49" "1.989.38"
Set-TextValue "E49" "  +3.09%  "

Set-TextValue "E50" "  +0.01%  "

Set-TextValue "D51" "4.63"
Set-TextValue "E51" "  +1.45%  "
